$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("detection_template_csv")

# Update J3 formula
$ws.Range("J3").Formula = '=$H3&IF($D3="raw",IF($E3<>""," ","")&$E3,"")&IF($D3="count"," count","")&", by "&IF($C3="TAC","TAC",$C3)&IF($D3="raw"," result","")'

# Clear J5 content (remove the "fix this for Ct" note text), keep formatting/style
$ws.Range("J5").ClearContents()
$ws.Rows.Item(5).EntireRow.AutoFit()

# Update the active selection to J3
$ws.Range("J3").Select()
